$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.782.25'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '3.407.64'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.409.10'
$ws.Range('E8').Value = '  +1.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.547'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.52%  '
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.122'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.432'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('D13').Value = '3.989.72'
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('E14').Value = '  -3.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000191'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Value = '63.753.02'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').Value = '3.408.88'
$ws.Range('E18').Value = '  +1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.531'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000119'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +23.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.87%  '
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.80'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.46'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('D39').Value = '2.974.00'
$ws.Range('E39').Value = '  +6.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0756'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('E43').Value = '  -5.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.758'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.61%  '
$ws.Range('E48').Value = '  +2.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.20'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +20.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.834'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.09%  '
